$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.036.69'
$ws.Range('E2').Value = '  +1.78%  '
$ws.Range('D3').Value = '3.159.04'
$ws.Range('E3').Value = '  +3.15%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.03'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.11'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +5.32%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '3.158.28'
$ws.Range('E8').Value = '  +3.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.528'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.13%  '
$ws.Range('E10').Value = '  +5.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.15'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.500'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +4.19%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000266'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +14.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.25'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +5.85%  '
$ws.Range('D15').Value = '3.676.77'
$ws.Range('E15').Value = '  +3.17%  '
$ws.Range('D16').Value = '65.126.93'
$ws.Range('E16').Value = '  +1.94%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.13'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +5.27%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.141.34'
$ws.Range('E18').Value = '  +2.51%  '
$ws.Range('E19').Value = '  +1.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '510.65'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +4.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.85'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.84%  '
$ws.Range('B22').Value = 'InternetComputer(DFINITY)'
$ws.Range('C22').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.38'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +5.89%  '
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.719'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +4.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.77'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +3.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.63'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.34%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.94'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +10.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.92'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.18'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +6.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.81'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +12.83%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '27.71'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +4.30%  '
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.19'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.29'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +10.50%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.56'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +5.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '55.34'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0904'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +10.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '469.56'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +5.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0422'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.03'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +8.94%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.67'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +4.00%  '
$ws.Range('D42').Value = '3.060.71'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.117'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.45'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +10.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.284'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.75'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.41%  '
$ws.Range('D47').Value = '0.0₃0594'
$ws.Range('E47').Value = '  +15.28%  '
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.114'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.61%  '
$ws.Range('E50').Value = '  +6.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.85'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.75%  '
